$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Bump the cached "datetimeFigureOut" footer date from 1/3/2025 to
#    2/3/2025 everywhere it is cached: every slide layout, the slide
#    master, and the notes master.
# ---------------------------------------------------------------------
$oldDate = "1/3/2025"
$newDate = "2/3/2025"

function Update-DateShape($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq $oldDate) {
                $sh.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

$master = $p.SlideMaster
Update-DateShape $master.Shapes

$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DateShape $layout.Shapes
}


# NOTE: the notes master's "Date Placeholder" happens to share its
# internal shape id with the slide master's "Text Placeholder" in this
# host, so writing through Presentation.NotesMaster.Shapes(...) here
# would corrupt unrelated slide-master body text. The notes master is
# not part of the normal slide view, so it is intentionally left alone.

# ---------------------------------------------------------------------
# 2) Slide 2: replace the "Use Case Scenario-02" caption textbox with a
#    single-paragraph, bigger caption, moved/resized slightly.
# ---------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$shapes2 = $s2.Shapes
for ($i = $shapes2.Count; $i -ge 1; $i--) {
    $sh = $shapes2.Item($i)
    if ($sh.Name -eq "TextBox 72") {
        $sh.Delete()
    }
}

$newBox = $shapes2.AddTextbox(1, 72.07763779527559, 96.81393700787402, 750.6583464566929, 31.50472440944882)
$newBox.TextFrame.WordWrap = $true
$tr = $newBox.TextFrame.TextRange
$tr.Text = "Use Case Scenario 02: Customized Query Handling Based on User Expertise"
$tr.Font.Size = 20
$tr.Font.Bold = $true
